$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.686774611473083
$ws.Range("B1").Value = 2.125629425048828
$ws.Range("C1").Value = 2.298197507858276
$ws.Range("D1").Value = 2.653293371200562
$ws.Range("E1").Value = 3.050495624542236
